## Fruta / hortaliza, semanal
## Insert one new week of data (3 quality rows: Especial, Primera, Segunda)
## for "Femacal de La Calera" / Frutilla at the top of the existing data
## block (new rows 405-407), pushing the previously existing rows 405-427
## down to 408-430.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 405, shifting everything
# below (old rows 405:427) down to 408:430.
$ws.Rows("405:407").Insert()

# Columns that are constant for every "Femacal de La Calera" / Frutilla /
# "Provincia de Melipilla" record in this block - fill them in for the
# three new rows at once.
$ws.Range("A405:A407").Value2 = 3
$ws.Range("B405:B407").Value2 = "Femacal de La Calera"
$ws.Range("C405:C407").Value2 = "Coquimbo"
$ws.Range("D405:D407").Value2 = 44931
$ws.Range("E405:E407").Value2 = 5
$ws.Range("F405:F407").Value2 = "Fruta"
$ws.Range("G405:G407").Value2 = 100101
$ws.Range("H405:H407").Value2 = "Berries"
$ws.Range("I405:I407").Value2 = 100112025
$ws.Range("J405:J407").Value2 = "Frutilla"
$ws.Range("K405:K407").Value2 = "Sin especificar"
$ws.Range("Q405:Q407").Value2 = "`$/bandeja 7 kilos"
$ws.Range("R405:R407").Value2 = "Provincia de Melipilla"
$ws.Range("T405:T407").Value2 = 7

# Row 405: Especial
$ws.Range("L405").Value2 = "Especial"
$ws.Range("M405").Value2 = 60
$ws.Range("N405").Value2 = 9000
$ws.Range("O405").Value2 = 9000
$ws.Range("P405").Value2 = 9000
$ws.Range("S405").Value2 = 1286

# Row 406: Primera
$ws.Range("L406").Value2 = "Primera"
$ws.Range("M406").Value2 = 65
$ws.Range("N406").Value2 = 7000
$ws.Range("O406").Value2 = 7000
$ws.Range("P406").Value2 = 7000
$ws.Range("S406").Value2 = 1000

# Row 407: Segunda
$ws.Range("L407").Value2 = "Segunda"
$ws.Range("M407").Value2 = 40
$ws.Range("N407").Value2 = 5000
$ws.Range("O407").Value2 = 5000
$ws.Range("P407").Value2 = 5000
$ws.Range("S407").Value2 = 714
